$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 96 - this shifts the existing rows 96-170
# down to 97-171, preserving all of their original values/formatting.
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with its (new) data.
$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96, 3).Value = "Ñuble"
$ws.Cells.Item(96, 4).Value = 44582
$ws.Cells.Item(96, 4).NumberFormat = $ws.Cells.Item(97, 4).NumberFormat
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = 100112017
$ws.Cells.Item(96, 7).Value = "Apio"
$ws.Cells.Item(96, 8).Value = "Americana (o)"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 60
$ws.Cells.Item(96, 11).Value = 8000
$ws.Cells.Item(96, 12).Value = 8500
$ws.Cells.Item(96, 13).Value = 8250
$ws.Cells.Item(96, 14).Value = "$/docena de matas"
$ws.Cells.Item(96, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(96, 16).Value = 1375
$ws.Cells.Item(96, 17).Value = 6
$ws.Cells.Item(96, 18).Value = "Hortaliza"
